# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of the existing header cell (AC1) onto the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the Wins/Losses/Ties values for every data row (2 through 50)
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 97   # AD = column 30
    $ws.Cells.Item($row, 31).Value = 65   # AE = column 31
    $ws.Cells.Item($row, 32).Value = 0    # AF = column 32
}
